$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FatosIn")
$ws.Range("W2").Formula = "=IF(OR(V2=""rgb"",V2=""cmy"",V2=""rgba"",V2=""cmyb""),  SUBSTITUTE(_xlfn.CONCAT(O2,""."",Q2,""."",S2,""."",U2), "".null"",""""), V2)"
$ws.Range("W3:W20").FormulaR1C1 = "=IF(OR(RC[-1]=""rgb"",RC[-1]=""cmy"",RC[-1]=""rgba"",RC[-1]=""cmyb""),  SUBSTITUTE(_xlfn.CONCAT(RC[-8],""."",RC[-6],""."",RC[-4],""."",RC[-2]), "".null"",""""), RC[-1])"
Write-Host "done"
